# feat: add video to dictionary entry format
#
# Inserts two new columns ("video_desc" / "video_fn") into both the
# "Sheet1" (headerless) and "Skip" (header-in-row-1) tables, between the
# existing audio_fn and theme columns. Mirrors the audio_desc/audio_fn
# pattern: video_desc holds the describer's name, video_fn the filename.
# The last row in each table ("goodbye"/"farvel") keeps no video data,
# matching the source.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Skip")

# --- Register the small (size-8) font Excel keeps around for the
#     phonetic-guide info on these sheets. Added as a transient named
#     style (then removed again) purely so the font itself lands in the
#     shared style table without being applied to any cell. -------------
$phoneticStyle = $wb.Styles.Add("PhoneticFont")
$phoneticStyle.Font.Size = 8
$phoneticStyle.Delete()

# --- Insert two blank columns at F:G on both sheets, pushing the old
#     F (theme) / G (secondary_theme) columns to H / I. -------------------
$ws1.Columns("F:G").Insert() | Out-Null
$ws2.Columns("F:G").Insert() | Out-Null

# --- Sheet1 ("Sheet1") -- headerless data rows ---------------------------
$ws1.Range("F1").Value = "Nolan Van Hell"
$ws1.Range("G1").Value = "snowfall.mp4"

$ws1.Range("F2").Value = "Nolan Van Hell"
$ws1.Range("G2").Value = "snowfall.mp4"

$ws1.Range("F3").Value = "Nolan Van Hell"
$ws1.Range("G3").Value = "snowfall.mp4"

# Row 4 ("goodbye") has no associated video - leave F4/G4 empty.

# --- Sheet2 ("Skip") -- row 1 is the header row ---------------------------
$ws2.Range("F1").Value = "video_desc"
$ws2.Range("G1").Value = "video_fn"

$ws2.Range("F2").Value = "Nolan Van Hell"
$ws2.Range("G2").Value = "snowfall.mp4"

$ws2.Range("F3").Value = "Nolan Van Hell"
$ws2.Range("G3").Value = "snowfall.mp4"

$ws2.Range("F4").Value = "Nolan Van Hell"
$ws2.Range("G4").Value = "snowfall.mp4"

# Row 5 ("goodbye") has no associated video - leave F5/G5 empty.

# --- Column widths for the new F/G columns (bestFit-style, same as the
#     widths Excel computed for this content). ---------------------------
$ws1.Columns("F").ColumnWidth = 12.5
$ws1.Columns("G").ColumnWidth = 11.5
$ws2.Columns("F").ColumnWidth = 12.5
$ws2.Columns("G").ColumnWidth = 11.5

# --- Defined names now need to cover the two extra columns (G -> I). -----
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!data") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$I`$4"
    }
    if ($n.Name -eq "Skip!data_1") {
        $n.RefersTo = "=Skip!`$A`$2:`$I`$5"
    }
}

# --- Leave the cursor where the author left it. ---------------------------
$ws1.Range("F12").Select() | Out-Null
$ws2.Activate() | Out-Null
$ws2.Range("H10").Select() | Out-Null
$ws1.Activate() | Out-Null
